$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLICMJDEZ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T180115Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJHMEUCIQCRhidtmD32BjqtO9WacLihY5B39yAlnbGf9yh%2F5ZdF6QIgLqEhfxqovsy8cPydN45uyq3wnrNzfr%2FRduKlOe%2BHig8q5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDC3%2BTQGiAaYmoQ1wkiq6A82pc%2Bbo7j52LiDNMD%2FbebbXY1Mr%2BtrtOs1K%2BMOHvEM39KfNctfQ0prhWjXqVmMy0oa02sw7JF1SGr2QImmiArYmo7YryD1hlLn1VWcLEzWUpK%2Ff4auRtp6%2BMUrMM5bu1gADS6KOAyzuqhhh2tApfHxZb1hkt0Tan9O%2FgYu0V7KpYMRA31iDEFYTE7ZCpm6rNJIAnz63Bgz2O0S5r9WcfGqLqzg8MerlyD%2BIxxOTqd39tOu9y0mjDLxoSpVP78AkPRsvbb%2BM0XOqusTHZMovp2UMnYnWiBHIB3c7VNV66WYKxwCyFXABhTrolSpX%2BnuTQJJs0of%2BcRJBqS7z0fYD5tAD7KW0VokAvQqz0K9FWEERxB6z4hHIGgBbBY9TV%2F%2BJYPcS50SgloVuaK6vVbsPXRC5akGeE%2FkYMdSF5GuX03wwzquwfvorlXy%2F3UHhpuLPdmfjS8JYrwqbHxcE4qPrCKVharmZLopp3%2FLDVDkdOMckkhpufC1fCmP9mbFrUq8H2ttsCXCozG5QbS3k7IbKSjJnCwqV7fgnq9jW46KOBvFJY4ZUCm9NdKJJ00mv9BXEeiQ%2FvSDonV5WDewwre61vwY6pQHqP9NP0Qk9jf3i75uG3FlMfzn9Fs0IHyPEX%2BYTWNbQv0Qtut7DTPemwgV83qBHMovCvNEoiTU80tLwako6FN7wacPGgDapDILUkeKw%2FqwJJjcPmi37q49UZg1GiE80dmpd%2FEgDeHCUhbcDYmj693E12Kch42qiCsJ9lYpFUa3DKOLigp%2BrWx7bNqE7IrKWRrEm2u8HzxGuBoX5L2noTgC2ZQDwYow%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2002%3A01%3A15%20GMT&X-Amz-Signature=dc5380605b03c4062d691d44e4e5c14cd66f923be45d4b6fd06f10e52fa4888b'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLICMJDEZ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T180115Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJHMEUCIQCRhidtmD32BjqtO9WacLihY5B39yAlnbGf9yh%2F5ZdF6QIgLqEhfxqovsy8cPydN45uyq3wnrNzfr%2FRduKlOe%2BHig8q5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDC3%2BTQGiAaYmoQ1wkiq6A82pc%2Bbo7j52LiDNMD%2FbebbXY1Mr%2BtrtOs1K%2BMOHvEM39KfNctfQ0prhWjXqVmMy0oa02sw7JF1SGr2QImmiArYmo7YryD1hlLn1VWcLEzWUpK%2Ff4auRtp6%2BMUrMM5bu1gADS6KOAyzuqhhh2tApfHxZb1hkt0Tan9O%2FgYu0V7KpYMRA31iDEFYTE7ZCpm6rNJIAnz63Bgz2O0S5r9WcfGqLqzg8MerlyD%2BIxxOTqd39tOu9y0mjDLxoSpVP78AkPRsvbb%2BM0XOqusTHZMovp2UMnYnWiBHIB3c7VNV66WYKxwCyFXABhTrolSpX%2BnuTQJJs0of%2BcRJBqS7z0fYD5tAD7KW0VokAvQqz0K9FWEERxB6z4hHIGgBbBY9TV%2F%2BJYPcS50SgloVuaK6vVbsPXRC5akGeE%2FkYMdSF5GuX03wwzquwfvorlXy%2F3UHhpuLPdmfjS8JYrwqbHxcE4qPrCKVharmZLopp3%2FLDVDkdOMckkhpufC1fCmP9mbFrUq8H2ttsCXCozG5QbS3k7IbKSjJnCwqV7fgnq9jW46KOBvFJY4ZUCm9NdKJJ00mv9BXEeiQ%2FvSDonV5WDewwre61vwY6pQHqP9NP0Qk9jf3i75uG3FlMfzn9Fs0IHyPEX%2BYTWNbQv0Qtut7DTPemwgV83qBHMovCvNEoiTU80tLwako6FN7wacPGgDapDILUkeKw%2FqwJJjcPmi37q49UZg1GiE80dmpd%2FEgDeHCUhbcDYmj693E12Kch42qiCsJ9lYpFUa3DKOLigp%2BrWx7bNqE7IrKWRrEm2u8HzxGuBoX5L2noTgC2ZQDwYow%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2002%3A01%3A15%20GMT&X-Amz-Signature=ba895cf268c46c01a9e102bdfaddf81ba2cde7eb643fa506db06b29e7c6b10e9'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLICMJDEZ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T180115Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJHMEUCIQCRhidtmD32BjqtO9WacLihY5B39yAlnbGf9yh%2F5ZdF6QIgLqEhfxqovsy8cPydN45uyq3wnrNzfr%2FRduKlOe%2BHig8q5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDC3%2BTQGiAaYmoQ1wkiq6A82pc%2Bbo7j52LiDNMD%2FbebbXY1Mr%2BtrtOs1K%2BMOHvEM39KfNctfQ0prhWjXqVmMy0oa02sw7JF1SGr2QImmiArYmo7YryD1hlLn1VWcLEzWUpK%2Ff4auRtp6%2BMUrMM5bu1gADS6KOAyzuqhhh2tApfHxZb1hkt0Tan9O%2FgYu0V7KpYMRA31iDEFYTE7ZCpm6rNJIAnz63Bgz2O0S5r9WcfGqLqzg8MerlyD%2BIxxOTqd39tOu9y0mjDLxoSpVP78AkPRsvbb%2BM0XOqusTHZMovp2UMnYnWiBHIB3c7VNV66WYKxwCyFXABhTrolSpX%2BnuTQJJs0of%2BcRJBqS7z0fYD5tAD7KW0VokAvQqz0K9FWEERxB6z4hHIGgBbBY9TV%2F%2BJYPcS50SgloVuaK6vVbsPXRC5akGeE%2FkYMdSF5GuX03wwzquwfvorlXy%2F3UHhpuLPdmfjS8JYrwqbHxcE4qPrCKVharmZLopp3%2FLDVDkdOMckkhpufC1fCmP9mbFrUq8H2ttsCXCozG5QbS3k7IbKSjJnCwqV7fgnq9jW46KOBvFJY4ZUCm9NdKJJ00mv9BXEeiQ%2FvSDonV5WDewwre61vwY6pQHqP9NP0Qk9jf3i75uG3FlMfzn9Fs0IHyPEX%2BYTWNbQv0Qtut7DTPemwgV83qBHMovCvNEoiTU80tLwako6FN7wacPGgDapDILUkeKw%2FqwJJjcPmi37q49UZg1GiE80dmpd%2FEgDeHCUhbcDYmj693E12Kch42qiCsJ9lYpFUa3DKOLigp%2BrWx7bNqE7IrKWRrEm2u8HzxGuBoX5L2noTgC2ZQDwYow%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2002%3A01%3A15%20GMT&X-Amz-Signature=fbc0b1363d41f6229373f8fe33090bea2d861786f2c4abb640c2ae37e97c5605'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLICMJDEZ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T180115Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJHMEUCIQCRhidtmD32BjqtO9WacLihY5B39yAlnbGf9yh%2F5ZdF6QIgLqEhfxqovsy8cPydN45uyq3wnrNzfr%2FRduKlOe%2BHig8q5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDC3%2BTQGiAaYmoQ1wkiq6A82pc%2Bbo7j52LiDNMD%2FbebbXY1Mr%2BtrtOs1K%2BMOHvEM39KfNctfQ0prhWjXqVmMy0oa02sw7JF1SGr2QImmiArYmo7YryD1hlLn1VWcLEzWUpK%2Ff4auRtp6%2BMUrMM5bu1gADS6KOAyzuqhhh2tApfHxZb1hkt0Tan9O%2FgYu0V7KpYMRA31iDEFYTE7ZCpm6rNJIAnz63Bgz2O0S5r9WcfGqLqzg8MerlyD%2BIxxOTqd39tOu9y0mjDLxoSpVP78AkPRsvbb%2BM0XOqusTHZMovp2UMnYnWiBHIB3c7VNV66WYKxwCyFXABhTrolSpX%2BnuTQJJs0of%2BcRJBqS7z0fYD5tAD7KW0VokAvQqz0K9FWEERxB6z4hHIGgBbBY9TV%2F%2BJYPcS50SgloVuaK6vVbsPXRC5akGeE%2FkYMdSF5GuX03wwzquwfvorlXy%2F3UHhpuLPdmfjS8JYrwqbHxcE4qPrCKVharmZLopp3%2FLDVDkdOMckkhpufC1fCmP9mbFrUq8H2ttsCXCozG5QbS3k7IbKSjJnCwqV7fgnq9jW46KOBvFJY4ZUCm9NdKJJ00mv9BXEeiQ%2FvSDonV5WDewwre61vwY6pQHqP9NP0Qk9jf3i75uG3FlMfzn9Fs0IHyPEX%2BYTWNbQv0Qtut7DTPemwgV83qBHMovCvNEoiTU80tLwako6FN7wacPGgDapDILUkeKw%2FqwJJjcPmi37q49UZg1GiE80dmpd%2FEgDeHCUhbcDYmj693E12Kch42qiCsJ9lYpFUa3DKOLigp%2BrWx7bNqE7IrKWRrEm2u8HzxGuBoX5L2noTgC2ZQDwYow%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2002%3A01%3A15%20GMT&X-Amz-Signature=e6f88459670772013f65a9cc0ddcb35eb851872155c210474dbb29ac82f6a19e'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLICMJDEZ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T180115Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJHMEUCIQCRhidtmD32BjqtO9WacLihY5B39yAlnbGf9yh%2F5ZdF6QIgLqEhfxqovsy8cPydN45uyq3wnrNzfr%2FRduKlOe%2BHig8q5gMI2%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDC3%2BTQGiAaYmoQ1wkiq6A82pc%2Bbo7j52LiDNMD%2FbebbXY1Mr%2BtrtOs1K%2BMOHvEM39KfNctfQ0prhWjXqVmMy0oa02sw7JF1SGr2QImmiArYmo7YryD1hlLn1VWcLEzWUpK%2Ff4auRtp6%2BMUrMM5bu1gADS6KOAyzuqhhh2tApfHxZb1hkt0Tan9O%2FgYu0V7KpYMRA31iDEFYTE7ZCpm6rNJIAnz63Bgz2O0S5r9WcfGqLqzg8MerlyD%2BIxxOTqd39tOu9y0mjDLxoSpVP78AkPRsvbb%2BM0XOqusTHZMovp2UMnYnWiBHIB3c7VNV66WYKxwCyFXABhTrolSpX%2BnuTQJJs0of%2BcRJBqS7z0fYD5tAD7KW0VokAvQqz0K9FWEERxB6z4hHIGgBbBY9TV%2F%2BJYPcS50SgloVuaK6vVbsPXRC5akGeE%2FkYMdSF5GuX03wwzquwfvorlXy%2F3UHhpuLPdmfjS8JYrwqbHxcE4qPrCKVharmZLopp3%2FLDVDkdOMckkhpufC1fCmP9mbFrUq8H2ttsCXCozG5QbS3k7IbKSjJnCwqV7fgnq9jW46KOBvFJY4ZUCm9NdKJJ00mv9BXEeiQ%2FvSDonV5WDewwre61vwY6pQHqP9NP0Qk9jf3i75uG3FlMfzn9Fs0IHyPEX%2BYTWNbQv0Qtut7DTPemwgV83qBHMovCvNEoiTU80tLwako6FN7wacPGgDapDILUkeKw%2FqwJJjcPmi37q49UZg1GiE80dmpd%2FEgDeHCUhbcDYmj693E12Kch42qiCsJ9lYpFUa3DKOLigp%2BrWx7bNqE7IrKWRrEm2u8HzxGuBoX5L2noTgC2ZQDwYow%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2002%3A01%3A15%20GMT&X-Amz-Signature=4d37bc9e19f1a08ac88bd8609e9d84603345e85a2b6dab30817b7167b124491c'
